# Add a new day (2025-11-05, serial 45966) of data to row 103 on every
# sheet, and fill in the previously-placeholder B102 value (2025-11-04)
# which had been uploaded as 0.
#
# Per-sheet (in tab order) values:
#   Sheet 1 (진양산업):   B102=3072  B103=3114
#   Sheet 2 (넥스트아이): B102=1169  B103=1159
#   Sheet 3 (삼보산업):   B102=1161  B103=1226
#   Sheet 4 (YBM넷):      B102=1851  B103=1869
#   Sheet 5 (NE능률):     B102=734   B103=734
#   Sheet 6 (위즈코프):   B102=1494  B103=1492
#   Sheet 7 (대영포장):   B102=2856  B103=2823

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"
$newDateSerial = 45966

$b102Values = @(3072, 1169, 1161, 1851, 734, 1494, 2856)
$b103Values = @(3114, 1159, 1226, 1869, 734, 1492, 2823)

for ($i = 1; $i -le 7; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Fill in the real value for the last existing row (2025-11-04), which
    # had been uploaded with a placeholder 0.
    $ws.Range("B102").Value = $b102Values[$i - 1]

    # Append the new row for 2025-11-05.
    $ws.Range("A103").NumberFormat = $dateFormat
    $ws.Range("A103").Value = $newDateSerial
    $ws.Range("B103").Value = $b103Values[$i - 1]
}
